$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "RXNO_DEF"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").VerticalAlignment = -4160
$ws.Range("F1").Borders.LineStyle = 1
$ws.Range("F2").Value = "['A generically dependent continuant that is about some thing. [IAO]']"
$ws.Range("F3").Value = "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']"
$ws.Range("F4").Value = "[]"
$ws.Range("F5").Value = "['A directive information entity that describes an intended process endpoint. When part of a plan specification the concretization is realized in a planned process in which the bearer tries to effect the world so that the process endpoint is achieved. [IAO]']"
$ws.Range("F6").Value = "[]"
$ws.Range("F7").Value = "[]"
$ws.Range("F8").Value = "[]"
$ws.Range("F9").Value = "[]"
$ws.Range("F10").Value = "[]"
$ws.Range("F11").Value = "[]"
$ws.Range("F12").Value = "[]"
$ws.Range("F13").Value = "[]"
$ws.Range("F14").Value = "[]"
$ws.Range("F15").Value = "[]"
$ws.Range("F16").Value = "[]"
$ws.Range("F17").Value = "[]"
$ws.Range("F18").Value = "[]"
